$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Season 14, matchdays prepares" - Оксанич Кирилл withdrew from the
# matchday sheet, so his whole row (the old row 15) is removed; every
# row below it (old rows 16-19) shifts up by one, and his now-unused
# name is dropped from the shared-string table automatically.
$ws.Rows.Item(15).Delete()

# Keep the worksheet's remembered sort range/condition in sync with the
# new (one-row-shorter) data range.
$sortObj = $ws.Sort
$sortFields = $sortObj.SortFields
$sortFields.Clear()
$sortFields.Add($ws.Range("B1:B18"))
$sortObj.SetRange($ws.Range("A1:W18"))
$sortObj.Header = 1
$sortObj.Apply()

# The author's last on-sheet selection ends up on the (new) row 15.
$ws.Range("A15:XFD15").Select()
